$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Import")

# Keep this sheet active/selected, as it already was in the source file.
$ws.Activate()

# Copy the style/formatting from the last existing data row (626) onto the
# new rows (627-634) so the new cells pick up style index 1 (wrap text),
# matching the rest of the translation table.
$ws.Range("A626:C626").Copy() | Out-Null
$ws.Range("A627:C634").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- New translation rows (language / key / Czech translation) -----------

# Row 627 - "glow" label (value entered before key)
$ws.Cells.Item(627, 1).Value = "cs"
$ws.Cells.Item(627, 3).Value = "Žhavení"
$ws.Cells.Item(627, 2).Value = "lab.build.glow.label"

# Row 628 - "glow" tooltip (tall row, long wrapped text)
$ws.Cells.Item(628, 1).Value = "cs"
$ws.Cells.Item(628, 2).Value = "lab.build.glow.label.tooltip"
$ws.Cells.Item(628, 3).Value = 'Tato hodnota udává rychlost žhavení; čím vyšší číslo, tím rychleji se spirálka rozžhaví; smyslem je poskytnout náhled, jak moc "divoký" build je.'

# Row 629 - "glow" table header (reuses the "Žhavení" translation)
$ws.Cells.Item(629, 1).Value = "cs"
$ws.Cells.Item(629, 2).Value = "lab.build.table.glow"
$ws.Cells.Item(629, 3).Value = "Žhavení"

# Row 630 - activate mixture button
$ws.Cells.Item(630, 1).Value = "cs"
$ws.Cells.Item(630, 2).Value = "lab.mixture.button.activate"
$ws.Cells.Item(630, 3).Value = "Aktivovat mix"

# Row 631 - deactivate mixture button
$ws.Cells.Item(631, 1).Value = "cs"
$ws.Cells.Item(631, 2).Value = "lab.mixture.button.deactivate"
$ws.Cells.Item(631, 3).Value = "Deaktivovat mix"

# Row 632 - deactivated success message
$ws.Cells.Item(632, 1).Value = "cs"
$ws.Cells.Item(632, 2).Value = "lab.mixture.deactivated.success"
$ws.Cells.Item(632, 3).Value = "Mix [{{data.name}}] byl úspěšně deaktivován."

# Row 633 - activated success message
$ws.Cells.Item(633, 1).Value = "cs"
$ws.Cells.Item(633, 2).Value = "lab.mixture.activated.success"
$ws.Cells.Item(633, 3).Value = "Mix [{{data.name}}] byl úspěšně aktivován."

# Row 634 - "active" preview flag
$ws.Cells.Item(634, 1).Value = "cs"
$ws.Cells.Item(634, 2).Value = "lab.mixture.preview.active"
$ws.Cells.Item(634, 3).Value = "Aktivní"

# Row 628 holds the long tooltip text and is taller (matches the source
# file's wrapped-text row height of 26.25 points).
$ws.Rows.Item(628).RowHeight = 26.25

# Restore view/selection state to the newly added area, similar to the
# source edit (scrolled down, new cell selected after appending rows).
$ws.Range("B629").Select() | Out-Null
